$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-04-29 Tuesday" "2025-04-30 Wednesday"
Replace-Text "769×2=1538" "511×8=4088"
Replace-Text "154×5=770" "104×4=416"
Replace-Text "208×5=1040" "738×2=1476"
Replace-Text "422×8=3376" "172×2=344"
Replace-Text "765×7=5355" "404×2=808"
Replace-Text "222×4=888" "479×8=3832"
Replace-Text "559×9=5031" "335×5=1675"
Replace-Text "269×6=1614" "101×9=909"
Replace-Text "807×8=6456" "287×8=2296"
Replace-Text "400×5=2000" "839×6=5034"
Replace-Text "413×2=826" "998×9=8982"
Replace-Text "407×8=3256" "940×7=6580"
Replace-Text "251×2=502" "710×8=5680"
Replace-Text "440×3=1320" "410×9=3690"
Replace-Text "727×4=2908" "540×9=4860"
Replace-Text "321×6=1926" "300×5=1500"
Replace-Text "224×7=1568" "113×5=565"
Replace-Text "393×4=1572" "700×8=5600"
Replace-Text "678×2=1356" "267×5=1335"
Replace-Text "181×4=724" "433×9=3897"
Replace-Text "349×8=2792" "847×7=5929"
Replace-Text "125×9=1125" "631×8=5048"
Replace-Text "194×3=582" "357×2=714"
Replace-Text "320×7=2240" "795×3=2385"
Replace-Text "980×5=4900" "709×6=4254"
